$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4): set the Value cell (B4) to the generated type name
$ws.Range("B4").Value = "SectiontableaucnopVs"

# "Date" row (row 8): update the generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
